# EscolherConfiguracaoOtima.xlsx — renumber the "Cenário Normal" use-case
# narrative: two steps are removed from the flow (the old "5.Indica quantia"
# actor step and the old "12. Confirma compra" system step), the
# pre-condition text is expanded, and every remaining step is renumbered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that disappear from the narrative -----------------
# Row 11 ("5.Indica quantia") goes away entirely; everything below shifts up.
$ws.Rows("11").Delete()
# After that shift, the old row 18 ("12. Confirma compra") is now row 17.
$ws.Rows("17").Delete()

# --- Pré condição: now also mentions the chosen-configuration flow ---------
$ws.Range("C4").Value = "Autenticado no sistema e vai escolher configuração ótima para comprar carro"

# --- Renumber / reword the remaining Actor input / System response steps ---
$ws.Range("D8").Value = "2.Pede quantia"
$ws.Range("C9").Value = "3.Indica quantia"
$ws.Range("D10").Value = "4. Regista quantia"
$ws.Range("D11").Value = "5. Calcula componentes usados na configuração ótima"
$ws.Range("D12").Value = "6. Calcula preço"
$ws.Range("D13").Value = "7. Mostra melhor carro, preço e suas especificações"
$ws.Range("D14").Value = "8. Pergunta se pretende confirmar"

# Row 15 used to hold the System response ("10. Pergunta..."); it now holds
# the Actor input "9. Confirma compra" instead, so move the value to column C.
$ws.Range("D15").ClearContents()
$ws.Range("C15").Value = "9. Confirma compra"

# Row 16 used to hold the Actor input ("11. Confirma compra"); it now holds
# the System response "10. Confirma compra" instead, so move it to column D.
$ws.Range("C16").ClearContents()
$ws.Range("D16").Value = "10. Confirma compra"

$ws.Range("D17").Value = "11. Insere carro na fila para produção"

# --- Exceção block (renumbered from "passo 11" to "passo 10") --------------
$ws.Range("B18").Value = "Exceção 1 [Não confirma compra] (passo 10)"
$ws.Range("C18").Value = "10.1 Não confirma compra"

# --- Refresh the on-screen selection to match the new row count ------------
$ws.Range("B18:B21").Select()

Write-Output "EscolherConfiguracaoOtima narrative updated"
